# Apply updated experiment values ("expermits todos no convexos menos el 5to")
$wb = $excel.ActiveWorkbook

$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")
$wsPunto    = $wb.Worksheets.Item("Punto_modificado")
# NOTE: "Vector_bf" and "Vector_BF" differ only by case, and Worksheets.Item(name)
# lookups are case-insensitive here, so both names would resolve to the same
# (first) sheet. Use the known 1-based sheet index instead to disambiguate.
$wsVecBf    = $wb.Worksheets.Item(5)   # Vector_bf
$wsVecBFup  = $wb.Worksheets.Item(6)   # Vector_BF

# Helper: write a value as TEXT (shared string), matching the original
# workbook where all these cells - even the numeric-looking ones - are
# stored as text (t="s") rather than numbers.
function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# Restricciones_del_follower (sheet3): rows 2-6, columns A,B,D,E,F updated; C unchanged
Set-TextValue $wsFollower.Range("A2") "4.49 - x - 0.5y"
Set-TextValue $wsFollower.Range("B2") "-2.49"
Set-TextValue $wsFollower.Range("D2") "0.62"
Set-TextValue $wsFollower.Range("E2") "6.8999999999999995"
Set-TextValue $wsFollower.Range("F2") "6.800000000000001"

Set-TextValue $wsFollower.Range("A3") "-4.4125 - 0.25x + y"
Set-TextValue $wsFollower.Range("B3") "2.4124999999999996"
Set-TextValue $wsFollower.Range("D3") "0.96"
Set-TextValue $wsFollower.Range("E3") "9.6"
Set-TextValue $wsFollower.Range("F3") "9.200000000000001"

Set-TextValue $wsFollower.Range("A4") "-4.49 + x + 0.5y"
Set-TextValue $wsFollower.Range("B4") "-3.51"
Set-TextValue $wsFollower.Range("D4") "0.88"
Set-TextValue $wsFollower.Range("E4") "4.1"
Set-TextValue $wsFollower.Range("F4") "4.699999999999999"

Set-TextValue $wsFollower.Range("A5") "-11.809999999999999 + x - 2y"
Set-TextValue $wsFollower.Range("B5") "-9.809999999999999"
Set-TextValue $wsFollower.Range("D5") "0.29"
Set-TextValue $wsFollower.Range("E5") "7.5"
Set-TextValue $wsFollower.Range("F5") "0.4"

Set-TextValue $wsFollower.Range("A6") "-4.92 - y"
Set-TextValue $wsFollower.Range("B6") "-4.92"
Set-TextValue $wsFollower.Range("D6") "0.34"
Set-TextValue $wsFollower.Range("E6") "5.2"
Set-TextValue $wsFollower.Range("F6") "7.800000000000001"

# Punto_modificado (sheet4)
Set-TextValue $wsPunto.Range("A2") "2.0300000000000002"
Set-TextValue $wsPunto.Range("B2") "4.92"

# Vector_bf (sheet5)
Set-TextValue $wsVecBf.Range("A2") "0.8300000000000001"

# Vector_BF (sheet6)
Set-TextValue $wsVecBFup.Range("A2") "-3.3000000000000007"
Set-TextValue $wsVecBFup.Range("A3") "11.0"
